$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "65.854.04"
$c.ClearFormats()
$ws.Range("E2").Value = "  +2.02%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.661.80"
$c.ClearFormats()
$ws.Range("E3").Value = "  +1.16%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "609.79"
$c.ClearFormats()
$ws.Range("E5").Value = "  +2.51%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "157.30"
$c.ClearFormats()
$ws.Range("E6").Value = "  +3.04%  "

$ws.Range("E7").Value = "  -0.07%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.589"
$c.ClearFormats()
$ws.Range("E8").Value = "  +0.03%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.660.22"
$c.ClearFormats()
$ws.Range("E9").Value = "  +1.20%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.124"
$c.ClearFormats()
$ws.Range("E10").Value = "  +8.21%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.405"
$c.ClearFormats()
$ws.Range("E11").Value = "  +1.97%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "5.91"
$c.ClearFormats()
$ws.Range("E12").Value = "  +1.95%  "

$ws.Range("E13").Value = "  +1.49%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "30.13"
$c.ClearFormats()
$ws.Range("E14").Value = "  +5.89%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0000200"
$c.ClearFormats()
$ws.Range("E15").Value = "  +16.17%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.141.45"
$c.ClearFormats()
$ws.Range("E16").Value = "  +1.27%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "65.593.46"
$c.ClearFormats()
$ws.Range("E17").Value = "  +1.77%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.667.93"
$c.ClearFormats()
$ws.Range("E18").Value = "  +0.54%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.78"
$c.ClearFormats()
$ws.Range("E19").Value = "  +4.03%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.92"
$c.ClearFormats()
$ws.Range("E20").Value = "  +2.91%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "360.43"
$c.ClearFormats()
$ws.Range("E21").Value = "  +2.99%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.46"
$c.ClearFormats()
$ws.Range("E22").Value = "  +4.79%  "

$ws.Range("E23").Value = "  +0.04%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "70.17"
$c.ClearFormats()
$ws.Range("E24").Value = "  +4.11%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.71"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.40%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "9.61"
$c.ClearFormats()
$ws.Range("E26").Value = "  +3.39%  "

$ws.Range("E27").Value = "  +16.38%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.65"
$c.ClearFormats()
$ws.Range("E28").Value = "  -0.52%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.168"
$c.ClearFormats()
$ws.Range("E29").Value = "  +2.77%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "8.19"
$c.ClearFormats()
$ws.Range("E30").Value = "  -0.69%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.23"
$c.ClearFormats()
$ws.Range("E31").Value = "  +6.40%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E32").Value = "  +0.12%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "531.12"
$c.ClearFormats()
$ws.Range("E33").Value = "  -3.18%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.82"
$c.ClearFormats()
$ws.Range("E34").Value = "  +0.66%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.60"
$c.ClearFormats()
$ws.Range("E35").Value = "  +0.83%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "6.41"
$c.ClearFormats()
$ws.Range("E36").Value = "  +3.45%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.434"
$c.ClearFormats()
$ws.Range("E37").Value = "  +2.35%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "20.72"
$c.ClearFormats()
$ws.Range("E38").Value = "  +3.05%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "163.06"
$c.ClearFormats()
$ws.Range("E39").Value = "  -0.64%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.01"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.32%  "

$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "166.35"
$c.ClearFormats()
$ws.Range("E43").Value = "  -1.42%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "42.03"
$c.ClearFormats()
$ws.Range("E44").Value = "  +0.89%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "4.16"
$c.ClearFormats()
$ws.Range("E45").Value = "  +1.54%  "

$ws.Range("E46").Value = "  +5.39%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0614"
$c.ClearFormats()
$ws.Range("E47").Value = "  +3.96%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "23.13"
$c.ClearFormats()
$ws.Range("E48").Value = "  -1.63%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0265"
$c.ClearFormats()
$ws.Range("E49").Value = "  +5.24%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.657"
$c.ClearFormats()
$ws.Range("E50").Value = "  +2.05%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0984"
$c.ClearFormats()
$ws.Range("E51").Value = "  +0.53%  "
